# Refresh the cryptos.xlsx price/volume snapshot (GitHub Actions scrape update).
#
# Column D ("Price") holds plain numeric-looking text (e.g. "113.05"); writing it
# straight through Range.Value would let Excel's smart-entry coerce it into a real
# number. To keep it as text (matching the source file, where every data cell is an
# unstyled string) we prefix the literal with an apostrophe - exactly how a user
# would force text entry in the UI - and then reset the cell's Style back to
# "Normal" so the transient quote-prefix style doesn't linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.684.15"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "2.293.73"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'113.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +16.81%  "
$ws.Range("D6").Value = "'269.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D9").Value = "'0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").Value = "'48.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.03%  "
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").Value = "'9.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.81%  "
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "'15.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "2.636.08"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "'0.852"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "2.300.17"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "43.701.11"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").Value = "'6.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.08%  "
$ws.Range("D21").Value = "'72.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "'2.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "'232.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "'9.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.87%  "
$ws.Range("D25").Value = "'2.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.43%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +3.49%  "
$ws.Range("D28").Value = "'42.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.37%  "
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "'2.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").Value = "'175.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'21.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("E34").Value = "  +4.80%  "
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "'4.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.87%  "
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("D38").Value = "'0.107"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").Value = "'3.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.48%  "
$ws.Range("B40").Value = "MultiversX"
$ws.Range("C40").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D40").Value = "'74.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.42%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "'13.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.42%  "
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").Value = "'2.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.47%  "
$ws.Range("D44").Value = "'6.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +22.42%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("D47").Value = "'8.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  +5.72%  "
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "'1.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("D51").Value = "'0.467"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.05%  "
